# Add season record columns (Wins, Losses, Ties) to the LAA_2013 worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - match the look of the existing header cells
# (bold font, thin border all around, centered horizontally, top vertically)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows 2-51: every team/row gets the same season record (78-84-0)
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 78   # AD
    $ws.Cells.Item($row, 31).Value = 84   # AE
    $ws.Cells.Item($row, 32).Value = 0    # AF
}
